$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7
$ws.Range("C7").Value = 7017
$ws.Range("E7").Value = 290920261

# Row 64
$ws.Range("C64").Value = 5219
$ws.Range("E64").Value = 20481621

# Row 91
$ws.Range("C91").Value = 151203
$ws.Range("E91").Value = 482928044

# Row 92
$ws.Range("C92").Value = 409295
$ws.Range("E92").Value = 1597281123

# Row 93
$ws.Range("C93").Value = 209658
$ws.Range("E93").Value = 1310059700

# Row 94
$ws.Range("C94").Value = 94235
$ws.Range("E94").Value = 919215690

# Row 95
$ws.Range("C95").Value = 50805
$ws.Range("E95").Value = 934508381

# Row 96
$ws.Range("C96").Value = 17323
$ws.Range("E96").Value = 797283923

# Row 97
$ws.Range("C97").Value = 2163
$ws.Range("E97").Value = 214451252

# Row 104
$ws.Range("C104").Value = 135301
$ws.Range("E104").Value = 272651977

# Row 128
$ws.Range("C128").Value = 25
$ws.Range("E128").Value = 509896

# Row 132
$ws.Range("C132").Value = 30292
$ws.Range("E132").Value = 174228945

# Row 135
$ws.Range("C135").Value = 1857
$ws.Range("E135").Value = 65742275

$wb.Save()
